$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.986.89"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "2.322.52"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "112.54"
$ws.Range("E5").Value = "  +18.14%  "
$ws.Range("D6").Value = "271.38"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E9").Value = "  +3.34%  "
$ws.Range("D10").Value = "47.23"
$ws.Range("E10").Value = "  +7.18%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "8.98"
$ws.Range("E12").Value = "  +15.72%  "
$ws.Range("D14").Value = "15.89"
$ws.Range("E14").Value = "  +4.47%  "
$ws.Range("D15").Value = "2.665.78"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "0.868"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "2.322.08"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "43.957.51"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "6.77"
$ws.Range("E20").Value = "  +9.16%  "
$ws.Range("D21").Value = "72.81"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +6.18%  "
$ws.Range("D23").Value = "235.21"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "9.58"
$ws.Range("E24").Value = "  +7.20%  "
$ws.Range("D25").Value = "2.88"
$ws.Range("E25").Value = "  +15.73%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "11.63"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").Value = "43.08"
$ws.Range("E28").Value = "  +14.69%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "178.14"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0949"
$ws.Range("E32").Value = "  +7.30%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "22.01"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  +8.50%  "
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  +21.81%  "
$ws.Range("D39").Value = "0.0361"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").Value = "0.248"
$ws.Range("E40").Value = "  +4.97%  "
$ws.Range("D41").Value = "2.42"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("D42").Value = "69.58"
$ws.Range("E42").Value = "  +12.02%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "1.40"
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "12.64"
$ws.Range("E45").Value = "  +5.97%  "
$ws.Range("D46").Value = "5.79"
$ws.Range("E46").Value = "  +10.71%  "
$ws.Range("D47").Value = "8.89"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "100.40"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").Value = "0.464"
$ws.Range("E51").Value = "  +9.45%  "
